$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update engine parameter values per "Add Inner and Outer Wall Curves" change
$ws.Range("B8").Value = 0.08889999999999999
$ws.Range("B9").Value = 0.2347036970157674
$ws.Range("B12").Value = 0.006207166618944346
$ws.Range("B13").Value = 0.1903706805793461
$ws.Range("B15").Value = 0.0443330164364213
